$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.515.03"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "2.637.49"
$ws.Range("E3").Value = "  +1.38%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.76%  "

$ws.Range("E10").Value = "  -1.16%  "

$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").Value = "3.104.50"
$ws.Range("E13").Value = "  +1.36%  "

$ws.Range("D14").Value = "59.419.23"
$ws.Range("E14").Value = "  +0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.66%  "

$ws.Range("D16").Value = "2.633.06"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("E18").Value = "  +2.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("E21").Value = "  -2.63%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.22%  "

$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("E35").Value = "  +2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.838"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.835"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "284.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("E44").Value = "  +2.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.85%  "

$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("E47").Value = "  +1.39%  "

$ws.Range("D48").Value = "1.958.52"
$ws.Range("E48").Value = "  +0.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
